# Daily tasks solving (Contest day)
# Add three new LeetCode entries to the "Main" sheet (rows 267-269),
# the same way previous rows were appended: type the values in, then
# re-assert column B's font so Excel registers a fresh cell style for
# the freshly pasted/typed title (mirrors the workbook's existing
# history of duplicated "Calibri" font/style entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Row 267: Number of Steps to Reduce a Number to Zero
$ws.Cells.Item(267, 1).Value = 1342
$ws.Cells.Item(267, 2).Value = "Number of Steps to Reduce a Number to Zero"
$ws.Cells.Item(267, 2).Font.Name = "Calibri"
$ws.Cells.Item(267, 3).Value = "Bit Monipulation"
$ws.Cells.Item(267, 4).Value = "Easy"
$ws.Cells.Item(267, 5).Value = "Solved"
$ws.Cells.Item(267, 6).Value = "Own"
$ws.Cells.Item(267, 7).Value = 43869

# Row 268: Number of Sub-arrays of Size K and Average Greater than or Equal to Threshold
$ws.Cells.Item(268, 1).Value = 1343
$ws.Cells.Item(268, 2).Value = "Number of Sub-arrays of Size K and Average Greater than or Equal to Threshold"
$ws.Cells.Item(268, 2).Font.Name = "Calibri"
$ws.Cells.Item(268, 3).Value = "Arrays"
$ws.Cells.Item(268, 4).Value = "Medium"
$ws.Cells.Item(268, 5).Value = "Solved"
$ws.Cells.Item(268, 6).Value = "Own"
$ws.Cells.Item(268, 7).Value = 43869

# Row 269: Angle Between Hands of a Clock
$ws.Cells.Item(269, 1).Value = 1344
$ws.Cells.Item(269, 2).Value = "Angle Between Hands of a Clock"
$ws.Cells.Item(269, 2).Font.Name = "Calibri"
$ws.Cells.Item(269, 3).Value = "Math"
$ws.Cells.Item(269, 4).Value = "Medium"
$ws.Cells.Item(269, 5).Value = "Solved"
$ws.Cells.Item(269, 6).Value = "Own"
$ws.Cells.Item(269, 7).Value = 43869

# Match date number format used by the rest of column G
$ws.Range("G267:G269").NumberFormat = "dd\.mm\.yyyy;@"

# Column B grew wider once the longest new title was added (matches the
# "best fit" width Excel computed for the new, longer title)
$ws.Columns.Item(2).ColumnWidth = 65.72135416666667

# Keep dimension/selection in sync with where editing ended
$ws.Activate()
[void]$ws.Range("B270").Select()
